# roman to int with math
$wb = $excel.ActiveWorkbook

# Work on the "数学" (Math) worksheet - this is the active sheet (4th tab)
$ws = $wb.Worksheets.Item("数学")

# New shared strings must be appended to the shared-string table in the same
# order the original author created them (F2, then E2, then C2, then D2) so
# that the resulting sharedStrings.xml indices line up with the target file.
$ws.Cells.Item(2, 6).Value = "O(n)，n是字符串长度"

$ws.Cells.Item(2, 5).Value = "哈希表`n字符串截取"

$ws.Cells.Item(2, 3).Value = "罗马数字包含以下七种字符: I， V， X， L，C，D 和 M。 `n字符          数值`nI             1`nV             5`nX             10`nL             50`nC             100`nD             500`nM             1000 `n例如， 罗马数字 2 写做 II ，即为两个并列的 1。12 写做 XII ，即为 X + II 。`n27 写做 XXVII, 即为 XX + V + I`n 。 `n通常情况下，罗马数字中小的数字在大的数字的右边。但也存在特例，例如 4 不写做 IIII，`n而是 IV。数字 1 在数字 5 的左边，所表示的数等于大数 5`n减小数 1 得到的数值 4 。同样地，数字 9 表示为 IX。这个特殊的规则只适用于以下六种情况： `nI 可以放在 V (5) 和 X (10) 的左边，来表示 4 和 9。 `nX 可以放在 L (50) 和 C (100) 的左边，来表示 40 和 90。 `nC 可以放在 D (500) 和 M (1000) 的左边，来表示 400 和 900。 `n给定一个罗马数字，将其转换成整数。输入确保在 1 到 3999 的范围内。 "

$ws.Cells.Item(2, 4).Value = "1 罗马字符与数字之间的映射关系以及六种特殊情况的映射关系`n2 从字符串头部先匹配两个字符串`n3 如果匹配成功就累加对应的数字，并且指针向前移动2;如果两个字符串匹配不成功，就匹配一个`n4 匹配一个如果成功就累加；不成功就跳过。`n5 注意如果字符串没有对应的罗马字符，要如何处理？在for循环中设定i++,无论是否执行for中的条件，for循环都会继续前进"

$ws.Cells.Item(2, 7).Value = "O(1)"

# The long, multi-line Roman-numeral text no longer fits at the old 22pt row
# height; with wrap text on, Excel would auto-grow the row to its maximum
# (409.6pt) to display it. Set that explicitly since content-based autosizing
# isn't available in this host.
$ws.Rows.Item(2).RowHeight = 409.6

# Update the selection to match the new active cell
$ws.Range("E2").Select() | Out-Null
